# Apply targeted odds updates to Sheet1, matching the author's commit
# "Atualizando o arquivo XLSX" (updating the odds/limits snapshot for 2025-12-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 2.28
$ws.Range("K2").Value = 3.7
$ws.Range("W2").Value = 1.31

# Row 3
$ws.Range("G3").Value = 2.02
$ws.Range("L3").Value = 1.36
$ws.Range("P3").Value = 1.8
$ws.Range("Q3").Value = 1.93
$ws.Range("AC3").Value = 990
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.59
$ws.Range("G5").Value = 1.72
$ws.Range("I5").Value = 8.199999999999999
$ws.Range("J5").Value = 3.8
$ws.Range("K5").Value = 4.4
$ws.Range("L5").Value = 1.37
$ws.Range("M5").Value = 1.07
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 1.91
$ws.Range("Q5").Value = 1.88
$ws.Range("R5").Value = 1.35
$ws.Range("S5").Value = 3.1
$ws.Range("T5").Value = 1.91
$ws.Range("U5").Value = 1.87
$ws.Range("V5").Value = 1.14
$ws.Range("W5").Value = 2.38
$ws.Range("X5").Value = 1000
$ws.Range("AB5").Value = 10.5
$ws.Range("AC5").Value = 11
$ws.Range("AG5").Value = 1000
$ws.Range("AJ5").Value = 1000

# Row 6
$ws.Range("H6").Value = 2.3
$ws.Range("I6").Value = 2.6
$ws.Range("J6").Value = 3.4
$ws.Range("N6").Value = 3.7
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 1.93
$ws.Range("Q6").Value = 1.89
$ws.Range("R6").Value = 1.36
$ws.Range("S6").Value = 3.2
$ws.Range("T6").Value = 1.7
$ws.Range("U6").Value = 2.16
$ws.Range("V6").Value = 1.62
$ws.Range("X6").Value = 18.5
$ws.Range("AA6").Value = 36
$ws.Range("AE6").Value = 28
$ws.Range("AI6").Value = 42
$ws.Range("AK6").Value = 40
$ws.Range("AL6").Value = 48
$ws.Range("AN6").Value = 38
$ws.Range("AO6").Value = 21

# Row 7
$ws.Range("F7").Value = 1.41

# Row 9
$ws.Range("N9").Value = 1.36
$ws.Range("P9").Value = 1.34
$ws.Range("T9").Value = 1.04
$ws.Range("U9").Value = 1.04
